$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'228.69"
$ws.Range("D3").Value = "'22.38"
$ws.Range("D4").Value = "'5.298"
$ws.Range("D5").Value = "'0.05529"
$ws.Range("D6").Value = "'3.391"
$ws.Range("D7").Value = "'6.464"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.053"
$ws.Range("E8").Value = "7FTXTokenFTT"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.7784"
$ws.Range("E9").Value = "8MXTokenMX"
$ws.Range("D10").Value = "'0.1378"
$ws.Range("D11").Value = "'0.07433"
$ws.Range("D12").Value = "'0.03143"
$ws.Range("D13").Value = "'0.02939"
$ws.Range("D14").Value = "'0.09263"
$ws.Range("D16").Value = "'3.260"
$ws.Range("D17").Value = "'0.04773"
$ws.Range("D18").Value = "'0.0005901"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006200"
$ws.Range("D20").Value = "'0.005227"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D23").Value = "'3.844"
$ws.Range("D24").Value = "'2.196"
$ws.Range("D27").Value = "'0.0005001"
$ws.Range("D40").Value = "'0.03951"
$ws.Range("D41").Value = "'0.007158"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003502"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1032"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.008571"
$ws.Range("D45").Value = "'0.00005444"
$ws.Range("D47").Value = "'0.7854"
$ws.Range("D48").Value = "'0.08802"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.01011"
